$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 271, shifting existing rows 271:317 down to 272:318
$ws.Rows.Item(271).Insert()

# Populate the newly inserted row 271 with the new record
$ws.Range("A271").Value = 5
$ws.Range("B271").Value = "Macroferia Regional de Talca"
$ws.Range("C271").Value = "Maule"
$ws.Range("D271").Value = 44694
$ws.Range("E271").Value = 7
$ws.Range("F271").Value = 100114014
$ws.Range("G271").Value = "Betarraga"
$ws.Range("H271").Value = "Sin especificar"
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 4000
$ws.Range("K271").Value = 600
$ws.Range("L271").Value = 600
$ws.Range("M271").Value = 600
$ws.Range("N271").Value = '$/paquete 5 unidades'
$ws.Range("O271").Value = "Región del Maule"
$ws.Range("P271").Value = 120
$ws.Range("Q271").Value = 5
$ws.Range("R271").Value = "Hortaliza"
